# Data update using git
# Apply the value changes to the "Resumo Inscrições Subsequente" worksheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E3").Value = 31

$ws.Range("E5").Value = 36
$ws.Range("F5").Value = 18
$ws.Range("H5").Value = 22

$ws.Range("E6").Value = 61

$ws.Range("E12").Value = 37

$ws.Range("F15").Value = 59
$ws.Range("H15").Value = 70

$ws.Range("E16").Value = 338
$ws.Range("F16").Value = 104
$ws.Range("H16").Value = 192

$ws.Range("E17").Value = 34
